$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-4 are cyclically rotated (row2 <- old row4, row3 <- old row2,
# row4 <- old row3), with the "Taxonsorteringsordning" (column B) receiving a
# fresh value in each destination row.

$ws.Range("A2").Value = 112182541
$ws.Range("B2").Value = 77402
$ws.Range("E2").Value = 6446
$ws.Range("F2").Value = "Kolflarnlav"
$ws.Range("G2").Value = "Carbonicola anthracophila"
$ws.Range("H2").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q2").Value = 756204

$ws.Range("A3").Value = 112181650
$ws.Range("B3").Value = 78713
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6458
$ws.Range("F3").Value = "Lunglav"
$ws.Range("G3").Value = "Lobaria pulmonaria"
$ws.Range("H3").Value = "(L.) Hoffm."
$ws.Range("Q3").Value = 756202
$ws.Range("R3").Value = 7291065

$ws.Range("A4").Value = 112181583
$ws.Range("B4").Value = 89794
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 65
$ws.Range("F4").Value = "Fläckporing"
$ws.Range("G4").Value = "Anthoporia albobrunnea"
$ws.Range("H4").Value = "(Romell) Karasiński & Niemelä"
$ws.Range("Q4").Value = 756188
$ws.Range("R4").Value = 7291007
